# Add a new "id" column (column D) with a sequential identifier for each
# data row, mirroring the row's position in the dataset (1-based).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 164

# Header for the new column.
$ws.Range("D1").Value = "id"

# Sequential ids for each data row (row 2 -> id 1, row 3 -> id 2, ...).
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 4).Value = ($r - 1)
}

# Restore view state that matches the edited workbook (active cell/selection
# and the scrolled position of the sheet).
$ws.Range("E163").Select()
